$d = $word.ActiveDocument

# Locate the three trailing paragraphs that must be removed:
#   1) the blank spacer paragraph right after "LOB1037: ..."
#   2) "Ver no Jupiter Salvar em pdf Salvar em docx"
#   3) the "(c) 2020 ... Creative Commons Attribution" paragraph
# by scanning paragraph text rather than relying on fixed indices.

$n = $d.Paragraphs.Count
$jupiterIdx = -1
$copyrightIdx = -1

for ($i = 1; $i -le $n; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*Ver no Jupiter*") {
        $jupiterIdx = $i
    }
    if ($t -like "*Powered by Jekyll*") {
        $copyrightIdx = $i
    }
}

if ($jupiterIdx -gt 0 -and $copyrightIdx -ge $jupiterIdx) {
    # Also swallow the blank paragraph immediately preceding the
    # "Ver no Jupiter ..." paragraph, so the "LOB1037" paragraph is
    # followed directly by whatever came after the copyright paragraph.
    $startPara = $jupiterIdx - 1
    if ($startPara -lt 1) {
        $startPara = $jupiterIdx
    }

    $delStart = $d.Paragraphs.Item($startPara).Range.Start
    $delEnd = $d.Paragraphs.Item($copyrightIdx).Range.End

    $r = $d.Range($delStart, $delEnd)
    $r.Delete()
}
